$wb = $excel.ActiveWorkbook

# --- ALC sheet: row 138 ---
# H138: 15283.8 -> 18105.25
# J138: 3998 -> 0
# L138: 11994 -> 0
# N138: removed entirely (cell deleted)
$wsAlc = $wb.Worksheets.Item("ALC")
$wsAlc.Range("H138").Value = 18105.25
$wsAlc.Range("J138").Value = 0
$wsAlc.Range("L138").Value = 0
$wsAlc.Range("N138").ClearContents()

# --- ARM sheet: rows 121-141 (except 136) ---
# Clear H:N for each row (removes the cells entirely)
$wsArm = $wb.Worksheets.Item("ARM")
$armRows = @(121,122,123,124,125,126,127,128,129,130,131,132,133,134,135,137,138,139,140,141)
foreach ($r in $armRows) {
    $wsArm.Range("H" + $r + ":N" + $r).ClearContents()
}

# --- WVR sheet: rows 119-141 (except 134) ---
# Clear H:N for each row (removes the cells entirely)
$wsWvr = $wb.Worksheets.Item("WVR")
$wvrRows = @(119,120,121,122,123,124,125,126,127,128,129,130,131,132,133,135,136,137,138,139,140,141)
foreach ($r in $wvrRows) {
    $wsWvr.Range("H" + $r + ":N" + $r).ClearContents()
}
